# Scheduled data refresh: update market-price driven profit columns (H-N)
# on Adamantoise_Profits across the per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 12134.182
$ws.Range("J17").Value = 13257.6
$ws.Range("L17").Value = 39772.8
$ws.Range("N17").Value = -40108.8
# Row 76
$ws.Range("H76").Value = 5969.8
$ws.Range("I76").Value = 5969.8
$ws.Range("K76").Value = 5969.8
$ws.Range("M76").Value = -5654.8
# Row 79
$ws.Range("H79").Value = 5969.8
$ws.Range("I79").Value = 5969.8
$ws.Range("K79").Value = 5969.8
$ws.Range("M79").Value = -4877.8
# Row 80
$ws.Range("H80").Value = 80011930
$ws.Range("I80").Value = 250000450
$ws.Range("J80").Value = 18197916
$ws.Range("K80").Value = 750001350
$ws.Range("L80").Value = 54593748
$ws.Range("M80").Value = -750000352
$ws.Range("N80").Value = -54595744
# Row 83
$ws.Range("H83").Value = 80011930
$ws.Range("I83").Value = 250000450
$ws.Range("J83").Value = 18197916
$ws.Range("K83").Value = 2250004050
$ws.Range("L83").Value = 163781244
$ws.Range("M83").Value = -2249999058
$ws.Range("N83").Value = -163791228
# Row 99
$ws.Range("H99").Value = 1257.8
$ws.Range("I99").Value = 799.36365
$ws.Range("J99").Value = 1818.1111
$ws.Range("K99").Value = 2398.09095
$ws.Range("L99").Value = 5454.3333
$ws.Range("M99").Value = -900.0909499999998
$ws.Range("N99").Value = -8450.3333
# Row 112
$ws.Range("H112").Value = 4417.88
$ws.Range("J112").Value = 4520.7085
$ws.Range("L112").Value = 13562.1255
$ws.Range("N112").Value = -15778.1255
# Row 120
$ws.Range("H120").Value = 202988.67
$ws.Range("J120").Value = 202988.67
$ws.Range("L120").Value = 202988.67
$ws.Range("N120").Value = -212664.67
# Row 128
$ws.Range("H128").Value = 75349.25
$ws.Range("J128").Value = 75349.25
$ws.Range("L128").Value = 75349.25
$ws.Range("N128").Value = -85309.25
# Row 132
$ws.Range("H132").Value = 4464.2974
$ws.Range("I132").Value = 4686.5806
$ws.Range("K132").Value = 14059.7418
$ws.Range("M132").Value = -11529.7418
# Row 136
$ws.Range("H136").Value = 89999
$ws.Range("I136").Value = 89999
$ws.Range("K136").Value = 89999
$ws.Range("M136").Value = -84899
# Row 138
$ws.Range("H138").Value = 4913.297
$ws.Range("J138").Value = 5576.5386
$ws.Range("L138").Value = 16729.6158
$ws.Range("N138").Value = -27009.6158
# Row 139
$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6503.93
$ws.Range("I32").Value = 5691.6665
$ws.Range("J32").Value = 25998.25
$ws.Range("K32").Value = 5691.6665
$ws.Range("L32").Value = 25998.25
$ws.Range("M32").Value = -5404.6665
$ws.Range("N32").Value = -26572.25
# Row 45
$ws.Range("H45").Value = 4884.4443
$ws.Range("I45").Value = 4601.4287
$ws.Range("K45").Value = 4601.4287
$ws.Range("M45").Value = -4224.4287
# Row 113
$ws.Range("H113").Value = 67677.60000000001
$ws.Range("J113").Value = 67677.60000000001
$ws.Range("L113").Value = 67677.60000000001
$ws.Range("N113").Value = -76355.60000000001
# Row 130
$ws.Range("H130").Value = 86816.71000000001
$ws.Range("J130").Value = 86816.71000000001
$ws.Range("L130").Value = 86816.71000000001
$ws.Range("N130").Value = -96856.71000000001
# Row 132
$ws.Range("H132").Value = 220921.73
$ws.Range("I132").Value = 273471.5
$ws.Range("K132").Value = 820414.5
$ws.Range("M132").Value = -817884.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2414.5833
$ws.Range("I105").Value = 1497.5
$ws.Range("K105").Value = 1497.5
$ws.Range("M105").Value = 249.5

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 525.5
$ws.Range("J22").Value = 751
$ws.Range("L22").Value = 751
$ws.Range("N22").Value = -1451
# Row 31
$ws.Range("H31").Value = 8984.697
$ws.Range("I31").Value = 6152.25
$ws.Range("K31").Value = 6152.25
$ws.Range("M31").Value = -5857.25
# Row 34
$ws.Range("H34").Value = 8984.697
$ws.Range("I34").Value = 6152.25
$ws.Range("K34").Value = 6152.25
$ws.Range("M34").Value = -5950.25
# Row 50
$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null
# Row 103
$ws.Range("H103").Value = 53792
$ws.Range("I103").Value = 13644.2
$ws.Range("K103").Value = 13644.2
$ws.Range("M103").Value = -12472.2
# Row 105
$ws.Range("H105").Value = 1991.6923
$ws.Range("I105").Value = 1626.6364
$ws.Range("K105").Value = 1626.6364
$ws.Range("M105").Value = 120.3635999999999
# Row 124
$ws.Range("H124").Value = 61130.332
$ws.Range("J124").Value = 61130.332
$ws.Range("L124").Value = 61130.332
$ws.Range("N124").Value = -66040.33199999999
# Row 130
$ws.Range("H130").Value = 49999.5
$ws.Range("J130").Value = 49999.5
$ws.Range("L130").Value = 49999.5
$ws.Range("N130").Value = -60039.5

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 565.6429000000001
$ws.Range("I7").Value = 619.9091
$ws.Range("J7").Value = 366.66666
$ws.Range("K7").Value = 1859.7273
$ws.Range("L7").Value = 1099.99998
$ws.Range("M7").Value = -1747.7273
$ws.Range("N7").Value = -1323.99998
# Row 17
$ws.Range("H17").Value = 1438.2307
$ws.Range("J17").Value = 684.6
$ws.Range("L17").Value = 2053.8
$ws.Range("N17").Value = -2391.8
# Row 68
$ws.Range("H68").Value = 1040.2
$ws.Range("J68").Value = 1065.6666
$ws.Range("L68").Value = 3196.9998
$ws.Range("N68").Value = -4818.9998
# Row 71
$ws.Range("H71").Value = 1040.2
$ws.Range("J71").Value = 1065.6666
$ws.Range("L71").Value = 9590.999400000001
$ws.Range("N71").Value = -17702.9994

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2874
$ws.Range("J126").Value = 2820.4443
$ws.Range("L126").Value = 8461.332900000001
$ws.Range("N126").Value = -13401.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 6749.5
$ws.Range("J16").Value = 6749.5
$ws.Range("L16").Value = 6749.5
$ws.Range("N16").Value = -7089.5
# Row 46
$ws.Range("H46").Value = 2537.5334
$ws.Range("J46").Value = 5164.1665
$ws.Range("L46").Value = 5164.1665
$ws.Range("N46").Value = -5540.1665
# Row 75
$ws.Range("H75").Value = 92558.25
$ws.Range("I75").Value = 45000
$ws.Range("J75").Value = 108411
$ws.Range("K75").Value = 45000
$ws.Range("L75").Value = 108411
$ws.Range("M75").Value = -44064
$ws.Range("N75").Value = -110283
# Row 78
$ws.Range("H78").Value = 92558.25
$ws.Range("I78").Value = 45000
$ws.Range("J78").Value = 108411
$ws.Range("K78").Value = 135000
$ws.Range("L78").Value = 325233
$ws.Range("M78").Value = -130320
$ws.Range("N78").Value = -334593
# Row 136
$ws.Range("H136").Value = 5628.8936
$ws.Range("I136").Value = 5261.6772
$ws.Range("J136").Value = 6340.375
$ws.Range("K136").Value = 15785.0316
$ws.Range("L136").Value = 19021.125
$ws.Range("M136").Value = -13235.0316
$ws.Range("N136").Value = -24121.125

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 69359
$ws.Range("J16").Value = 69359
$ws.Range("L16").Value = 69359
$ws.Range("N16").Value = -69943
# Row 68
$ws.Range("H68").Value = 74989
$ws.Range("J68").Value = 74989
$ws.Range("L68").Value = 74989
$ws.Range("N68").Value = -76611
# Row 71
$ws.Range("H71").Value = 74989
$ws.Range("J71").Value = 74989
$ws.Range("L71").Value = 224967
$ws.Range("N71").Value = -233079
# Row 113
$ws.Range("H113").Value = 339.2
$ws.Range("I113").Value = 202
$ws.Range("J113").Value = 373.5
$ws.Range("K113").Value = 606
$ws.Range("L113").Value = 1120.5
$ws.Range("M113").Value = 1564
$ws.Range("N113").Value = -5460.5
# Row 126
$ws.Range("H126").Value = 3848.8125
$ws.Range("I126").Value = 1824.2727
$ws.Range("K126").Value = 5472.8181
$ws.Range("M126").Value = -3002.8181
# Row 132
$ws.Range("H132").Value = 18664.791
$ws.Range("I132").Value = 25992.023
$ws.Range("K132").Value = 77976.069
$ws.Range("M132").Value = -75446.069
# Row 136
$ws.Range("H136").Value = 47546.477
$ws.Range("I136").Value = 3006.3076
$ws.Range("K136").Value = 9018.9228
$ws.Range("M136").Value = -6468.9228
